# Data Driven Extended With Docker
# Rename Sheet1 -> TestCases, Sheet2 -> TestData, remove Sheet3,
# populate both sheets with test-case / test-data tables, and restore
# the view/selection state captured in the target workbook.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item(1)
$wsData  = $wb.Worksheets.Item(2)

# --- rename sheets -------------------------------------------------------
$wsCases.Name = "TestCases"
$wsData.Name  = "TestData"

# --- drop the unused third sheet -----------------------------------------
$wb.Worksheets.Item(3).Delete()

# --- TestCases sheet content ---------------------------------------------
$wsCases.Range("A1").Value = "TestCases"
$wsCases.Range("B1").Value = "Runmode"
$wsCases.Range("A2").Value = "AddCustomerTest"
$wsCases.Range("B2").Value = "Y"
$wsCases.Range("A3").Value = "OpenAccountTest"
$wsCases.Range("B3").Value = "Y"

# column A autosize-ish manual width (character units, quantised by host)
$wsCases.Columns.Item(1).ColumnWidth = 18.833333333333332

# --- TestData sheet content ------------------------------------------------
$wsData.Range("A1").Value = "AddCustomerTest"

$wsData.Range("A2").Value = "Runmode"
$wsData.Range("B2").Value = "firstname"
$wsData.Range("C2").Value = "lastname"
$wsData.Range("D2").Value = "postcode"

$wsData.Range("A3").Value = "Y"
$wsData.Range("B3").Value = "Rahul"
$wsData.Range("C3").Value = "Dash"
$wsData.Range("D3").Value = 7878

$wsData.Range("A4").Value = "Y"
$wsData.Range("B4").Value = "Amit"
$wsData.Range("C4").Value = "jena"
$wsData.Range("D4").Value = 6787

$wsData.Range("A7").Value = "OpenAccountTest"

$wsData.Range("A8").Value = "Runmode"
$wsData.Range("B8").Value = "customer"
$wsData.Range("C8").Value = "currency"

$wsData.Range("A9").Value = "N"
$wsData.Range("B9").Value = "Rahul Dash"
$wsData.Range("C9").Value = "Rupee"

$wsData.Range("A10").Value = "Y"
$wsData.Range("B10").Value = "Amit Jena"
$wsData.Range("C10").Value = "Dollar"

$wsData.Range("A11").Value = "Y"
$wsData.Range("B11").Value = "Amit Jena"
$wsData.Range("C11").Value = "Dollar"

$wsData.Columns.Item(1).ColumnWidth = 16
$wsData.Columns.Item(2).ColumnWidth = 10.166666666666666

# --- selections & active sheet/tab ---------------------------------------
$wsCases.Range("B3").Select()
$wsData.Activate()
$wsData.Range("A10").Select()

# --- restore the saved window geometry (best effort; host may not persist
#     these into bookViews on export, but set them anyway for fidelity) ---
$win = $excel.ActiveWindow
$win.Left   = 0
$win.Top    = 330
$win.Width  = 15210
$win.Height = 2745
